# "End of Day (2)" — update the data-collection log (Sheet1) with the
# day's progress: flip a batch of per-year boolean flags from FALSE to
# TRUE (and one back to FALSE) across existing rows, and fill in the
# progress flags for the newly-tracked tickers in rows 98-107
# (CMCSA, SFIX, CHWY, RNG, ETSY, PRTS, DASH, UBER, U, LYFT) which had
# previously been blank placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Flags($row, $cols) {
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = [bool]$cols[$col]
    }
}

# Existing rows: flip specific year-flag cells (mostly 0 -> 1).
Set-Flags 3  @{ F = 1 }
Set-Flags 10 @{ D = 1; E = 1; F = 1 }
Set-Flags 43 @{ F = 1 }
Set-Flags 50 @{ E = 1; F = 1; G = 1; H = 1; I = 1; J = 1; K = 1 }
Set-Flags 58 @{ B = 1; C = 1; D = 1; E = 1; F = 1; G = 1; H = 1; I = 1; J = 1; K = 1 }
Set-Flags 62 @{ E = 1; F = 1; G = 1; H = 1; I = 1; J = 1 }
Set-Flags 72 @{ E = 1; F = 1; G = 1; H = 1 }
Set-Flags 83 @{ D = 1; E = 1; F = 1; G = 1; H = 1; J = 1; K = 1 }
Set-Flags 90 @{ M = 0 }
Set-Flags 95 @{ B = 1 }
Set-Flags 96 @{ I = 1; J = 1; K = 1 }

# Newly progressed rows 98-107 (CMCSA through LYFT) - fill in the
# year/flag columns that now have data collected.
Set-Flags 98  @{ F = 1; G = 1; H = 1; I = 1; J = 1; K = 1 }
Set-Flags 99  @{ B = 1; C = 1; D = 1; M = 1 }
Set-Flags 100 @{ B = 1; C = 1; M = 1 }
Set-Flags 101 @{ B = 1; C = 1; D = 1; E = 1; F = 1; G = 1; H = 1; M = 1 }
Set-Flags 102 @{ B = 1; C = 1; D = 1; E = 1; F = 1; G = 1; M = 1 }
Set-Flags 103 @{ B = 1; C = 1; D = 1; E = 1; F = 1; G = 1; H = 1; M = 1 }
Set-Flags 104 @{ B = 1; M = 1 }
Set-Flags 105 @{ B = 1; C = 1; M = 1 }
Set-Flags 106 @{ B = 1; M = 1 }
Set-Flags 107 @{ B = 1; C = 1; M = 1 }
